$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "595.57")
# are stored as text, matching the source data (t="inlineStr" in the original).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.710.78'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '3.834.35'
$ws.Range('E3').Value = '  -2.23%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '595.57'
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').Value = '165.63'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').Value = '3.834.60'
$ws.Range('E7').Value = '  -2.15%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.63%  '
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('D11').Value = '6.26'
$ws.Range('E11').Value = '  -2.40%  '
$ws.Range('E12').Value = '  -1.78%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').Value = '36.67'
$ws.Range('E14').Value = '  -1.70%  '
$ws.Range('D15').Value = '4.476.65'
$ws.Range('E15').Value = '  -2.15%  '
$ws.Range('D16').Value = '3.827.90'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').Value = '67.709.71'
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '18.24'
$ws.Range('E18').Value = '  +6.86%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = '7.40'
$ws.Range('E19').Value = '  -1.19%  '
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').Value = '10.64'
$ws.Range('E21').Value = '  -4.56%  '
$ws.Range('D22').Value = '465.93'
$ws.Range('E22').Value = '  -4.58%  '
$ws.Range('D23').Value = '0.725'
$ws.Range('E23').Value = '  +0.28%  '
$ws.Range('E24').Value = '  -4.83%  '
$ws.Range('E25').Value = '  -1.47%  '
$ws.Range('E26').Value = '  -3.49%  '
$ws.Range('D27').Value = '12.04'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '9.95'
$ws.Range('E29').Value = '  -2.17%  '
$ws.Range('D30').Value = '2.90'
$ws.Range('E30').Value = '  -1.74%  '
$ws.Range('D31').Value = '3.981.68'
$ws.Range('E31').Value = '  -2.28%  '
$ws.Range('D32').Value = '7.67'
$ws.Range('E32').Value = '  -2.44%  '
$ws.Range('E33').Value = '  -5.05%  '
$ws.Range('E34').Value = '  -4.54%  '
$ws.Range('D35').Value = '3.804.73'
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('E38').Value = '  -3.15%  '
$ws.Range('D39').Value = '5.85'
$ws.Range('E39').Value = '  -1.28%  '
$ws.Range('D40').Value = '3.23'
$ws.Range('E40').Value = '  +7.52%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  -3.77%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = '421.96'
$ws.Range('E43').Value = '  -3.61%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '1.96'
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('D46').Value = '47.18'
$ws.Range('E46').Value = '  -2.68%  '
$ws.Range('D47').Value = '8.49'
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').Value = '142.77'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0353'
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').Value = '0.000264'
$ws.Range('E50').Value = '  +10.81%  '
$ws.Range('D51').Value = '38.94'
$ws.Range('E51').Value = '  -1.06%  '

# Restore default styling on column D (remove the temporary text number format)
$ws.Range("D2:D51").Style = "Normal"
